$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# row hunk 0
$ws.Range("H6").Value = 9726.700000000001
$ws.Range("I6").Value = 10163.167
$ws.Range("J6").Value = 5798.5
$ws.Range("K6").Value = 30489.501
$ws.Range("L6").Value = 17395.5
$ws.Range("M6").Value = -30377.501
$ws.Range("N6").Value = -17619.5
# row hunk 1
$ws.Range("H17").Value = 3116164
$ws.Range("J17").Value = 3489875.8
$ws.Range("L17").Value = 10469627.4
$ws.Range("N17").Value = -10469963.4
# row hunk 2
$ws.Range("H40").Value = 2633.3333
$ws.Range("J40").Value = 2688.2354
$ws.Range("L40").Value = 2688.2354
$ws.Range("N40").Value = -3038.2354
# row hunk 3
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
# row hunk 4
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
# row hunk 5
$ws.Range("H98").Value = 2999
$ws.Range("I98").Value = 2999
$ws.Range("K98").Value = 2999
$ws.Range("M98").Value = -1501
# row hunk 6
$ws.Range("H122").Value = 2999
$ws.Range("I122").Value = 2999
$ws.Range("K122").Value = 8997
$ws.Range("M122").Value = -6547
# row hunk 7
$ws.Range("H132").Value = 3269.9062
$ws.Range("I132").Value = 3117.3225
$ws.Range("K132").Value = 9351.967500000001
$ws.Range("M132").Value = -6821.967500000001
# row hunk 8
$ws.Range("H137").Value = 18995.5
$ws.Range("I137").Value = 8749.799999999999
$ws.Range("K137").Value = 26249.4
$ws.Range("M137").Value = -23699.4

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# row hunk 9
$ws.Range("H32").Value = 1989.3771
$ws.Range("I32").Value = 1258.2909
$ws.Range("K32").Value = 1258.2909
$ws.Range("M32").Value = -971.2909
# row hunk 10
$ws.Range("H45").Value = 13989.667
$ws.Range("I45").Value = 15488.375
$ws.Range("J45").Value = 2000
$ws.Range("K45").Value = 15488.375
$ws.Range("L45").Value = 2000
$ws.Range("M45").Value = -15111.375
$ws.Range("N45").Value = -2754
# row hunk 11
$ws.Range("H63").Value = 2818.0967
$ws.Range("I63").Value = 2101.7856
$ws.Range("J63").Value = 9503.666999999999
$ws.Range("K63").Value = 2101.7856
$ws.Range("L63").Value = 9503.666999999999
$ws.Range("M63").Value = -1415.7856
$ws.Range("N63").Value = -10875.667
# row hunk 12
$ws.Range("H66").Value = 2818.0967
$ws.Range("I66").Value = 2101.7856
$ws.Range("J66").Value = 9503.666999999999
$ws.Range("K66").Value = 10508.928
$ws.Range("L66").Value = 47518.335
$ws.Range("M66").Value = -7076.928
$ws.Range("N66").Value = -54382.335
# row hunk 13
$ws.Range("H132").Value = 12134.4
$ws.Range("I132").Value = 10712.733
$ws.Range("J132").Value = 16399.4
$ws.Range("K132").Value = 32138.199
$ws.Range("L132").Value = 49198.2
$ws.Range("M132").Value = -29608.199
$ws.Range("N132").Value = -54258.2

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# row hunk 14
$ws.Range("H20").Value = 4091.2727
$ws.Range("I20").Value = 2395
$ws.Range("K20").Value = 2395
$ws.Range("M20").Value = -2148
# row hunk 15
$ws.Range("H105").Value = 2908.484
$ws.Range("I105").Value = 3068.5264
$ws.Range("K105").Value = 3068.5264
$ws.Range("M105").Value = -1321.5264
# row hunk 16
$ws.Range("H135").Value = 41122.473
$ws.Range("J135").Value = 41122.473
$ws.Range("L135").Value = 41122.473
$ws.Range("N135").Value = -51262.473
# row hunk 17
$ws.Range("H137").Value = 69998
$ws.Range("J137").Value = 69998
$ws.Range("L137").Value = 69998
$ws.Range("N137").Value = -80198
# row hunk 18
$ws.Range("H138").Value = 65568.71000000001
$ws.Range("J138").Value = 65568.71000000001
$ws.Range("L138").Value = 65568.71000000001
$ws.Range("N138").Value = -75848.71000000001
# row hunk 19
$ws.Range("H140").Value = 171614.39
$ws.Range("J140").Value = 171614.39
$ws.Range("L140").Value = 171614.39
$ws.Range("N140").Value = -181974.39

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# row hunk 20
$ws.Range("H2").Value = 43.333332
$ws.Range("I2").Value = 43.333332
$ws.Range("K2").Value = 43.333332
$ws.Range("M2").Value = 69.666668
# row hunk 21
$ws.Range("H5").Value = 2381.625
$ws.Range("I5").Value = 1310.6
$ws.Range("K5").Value = 1310.6
$ws.Range("M5").Value = -1198.6
# row hunk 22
$ws.Range("H134").Value = 10196.863
$ws.Range("I134").Value = 8205.6
$ws.Range("K134").Value = 24616.8
$ws.Range("M134").Value = -22081.8

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# row hunk 23
$ws.Range("H5").Value = 1536.2593
$ws.Range("I5").Value = 476.53845
$ws.Range("J5").Value = 2520.2856
$ws.Range("K5").Value = 1429.61535
$ws.Range("L5").Value = 7560.8568
$ws.Range("M5").Value = -1317.61535
$ws.Range("N5").Value = -7784.8568
# row hunk 24
$ws.Range("H116").Value = 1639.75
$ws.Range("J116").Value = 1840
$ws.Range("L116").Value = 5520
$ws.Range("N116").Value = -12404
# row hunk 25
$ws.Range("H135").Value = 1536.2593
$ws.Range("I135").Value = 476.53845
$ws.Range("J135").Value = 2520.2856
$ws.Range("K135").Value = 4288.84605
$ws.Range("L135").Value = 22682.5704
$ws.Range("M135").Value = -1753.84605
$ws.Range("N135").Value = -27752.5704

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# row hunk 26
$ws.Range("H102").Value = 1658.65
$ws.Range("I102").Value = 1534.3572
$ws.Range("K102").Value = 1534.3572
$ws.Range("M102").Value = 87.64280000000008
# row hunk 27
$ws.Range("H113").Value = 2987.125
$ws.Range("I113").Value = 2982.8333
$ws.Range("K113").Value = 2982.8333
$ws.Range("M113").Value = -812.8332999999998
# row hunk 28
$ws.Range("H132").Value = 3800.4
$ws.Range("I132").Value = 3922.7693
$ws.Range("J132").Value = 3005
$ws.Range("K132").Value = 11768.3079
$ws.Range("L132").Value = 9015
$ws.Range("M132").Value = -9238.3079
$ws.Range("N132").Value = -14075

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# row hunk 29
$ws.Range("H22").Value = 3083.2068
$ws.Range("I22").Value = 2285.1428
$ws.Range("K22").Value = 2285.1428
$ws.Range("M22").Value = -1990.1428
# row hunk 30
$ws.Range("H27").Value = 3083.2068
$ws.Range("I27").Value = 2285.1428
$ws.Range("K27").Value = 2285.1428
$ws.Range("M27").Value = -2178.1428
# row hunk 31
$ws.Range("H40").Value = 1573.1428
$ws.Range("I40").Value = 1335.4166
$ws.Range("K40").Value = 1335.4166
$ws.Range("M40").Value = -1199.4166
# row hunk 32
$ws.Range("H46").Value = 1436.7916
$ws.Range("I46").Value = 1035.6364
$ws.Range("J46").Value = 1776.2307
$ws.Range("K46").Value = 1035.6364
$ws.Range("L46").Value = 1776.2307
$ws.Range("M46").Value = -847.6364000000001
$ws.Range("N46").Value = -2152.2307
# row hunk 33
$ws.Range("H68").Value = 2903.0833
$ws.Range("I68").Value = 2589.3684
$ws.Range("K68").Value = 2589.3684
$ws.Range("M68").Value = -1840.3684
# row hunk 34
$ws.Range("H71").Value = 2903.0833
$ws.Range("I71").Value = 2589.3684
$ws.Range("K71").Value = 12946.842
$ws.Range("M71").Value = -9202.841999999999
# row hunk 35
$ws.Range("H132").Value = 10000
$ws.Range("I132").Value = 10000
$ws.Range("K132").Value = 30000
$ws.Range("M132").Value = -27470
# row hunk 36
$ws.Range("H136").Value = 4557.1
$ws.Range("I136").Value = 4012.2
$ws.Range("J136").Value = 6191.8
$ws.Range("K136").Value = 12036.6
$ws.Range("L136").Value = 18575.4
$ws.Range("M136").Value = -9486.599999999999
$ws.Range("N136").Value = -23675.4

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# row hunk 37
$ws.Range("H61").Value = 34898.168
$ws.Range("I61").Value = 34898.168
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 34898.168
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -34606.168
$ws.Range("N61").ClearContents()
# row hunk 38
$ws.Range("H100").Value = 1093.1666
$ws.Range("I100").Value = 1061.8
$ws.Range("K100").Value = 2123.6
$ws.Range("M100").Value = -1582.6
# row hunk 39
$ws.Range("H132").Value = 26013.645
$ws.Range("I132").Value = 18665.055
$ws.Range("K132").Value = 55995.165
$ws.Range("M132").Value = -53465.165
